$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.061.07"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.006.06"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'593.61"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "'147.15"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D8").Value = "3.005.12"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E10").Value = "  +6.78%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "3.500.13"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "62.019.48"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "3.005.45"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'446.38"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "'14.10"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "'0.686"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'82.18"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'11.08"
$ws.Range("E25").Value = "  +10.41%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'7.23"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "'27.38"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'50.14"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "'8.99"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "'2.94"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("D43").Value = "'41.49"
$ws.Range("E43").Value = "  +10.83%  "
$ws.Range("D44").Value = "'0.280"
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("D45").Value = "'392.43"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "2.716.16"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'133.59"
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("D50").Value = "'2.16"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("E51").Value = "  -1.61%  "
